$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank row 32 (A32:B32), shifting all rows below it up by one.
$ws.Rows.Item(32).Delete()

# Row 31 (unaffected by the row-32 deletion, since 31 < 32) is given an
# explicit custom height, matching the taller wrapped text.
$ws.Rows.Item(31).RowHeight = 36

# Update the view to match the post-edit state: scrolled so row 30 is near the
# top, with the row that is now 32 (previously row 33) selected as an entire
# row (as if the user clicked the row header).
$excel.ActiveWindow.ScrollRow = 30
$ws.Range("A32:XFD32").Select()
